# Car_Circuitry_BOM.xlsx edit script
# - Fix L1 inductor row (footprint / libref were duplicated from L2)
# - Add "Total" (G) and "Boards" (H) columns: Total = Quantity * Boards, Boards = 5 (entered once in H2)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fix the messed-up L1 inductor row (row 17) ---
# Footprint (D17) and LibRef (E17) had incorrectly been copied from the L2 row.
# Leading apostrophe preserves the existing quote-prefixed text style of these cells
# (otherwise the numeric-looking LibRef would be coerced into a plain number cell).
$ws.Range("D17").Value = "'INDPM6664X610N"
$ws.Range("E17").Value = "'74439346047"

# --- Add new "Total" / "Boards" columns ---
$ws.Range("G1").Value = "Total"
$ws.Range("H1").Value = "Boards"

# Number of boards being built (entered once, referenced by every row)
$ws.Range("H2").Value = 5

# Total = Quantity * $H$2 for every data row (2-29)
$ws.Range("G2").Formula = "=F2*`$H`$2"
$ws.Range("G3:G29").Formula = "=F3*`$H`$2"

# Match the formatting (borders/fill) used by the existing neighbouring columns
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)

$ws.Range("F2").Copy()
$ws.Range("G2:H2").PasteSpecial(-4122)

$ws.Range("F3").Copy()
$ws.Range("G3:G29").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Keep the active selection where Excel left it after the edit
$ws.Range("F4").Select()
